# Daily attendance processing - 2025-12-05 14:28:29
# Normalizes the "Recorded By" (column G) entries so that the leading
# author name is rotated to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System";
    "backup@backdoor.com, System" = "System, backup@backdoor.com";
    "backup@backdoor.com, system, System" = "system, System, backup@backdoor.com";
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
